$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2834
$ws.Range("I51").Value = 2500
$ws.Range("J51").Value = 3502
$ws.Range("K51").Value = 2500
$ws.Range("L51").Value = 3502
$ws.Range("M51").Value = -2016
$ws.Range("N51").Value = -4470

$ws.Range("H76").Value = 7209.9
$ws.Range("I76").Value = 6419.8
$ws.Range("K76").Value = 6419.8
$ws.Range("M76").Value = -6104.8

$ws.Range("H79").Value = 7209.9
$ws.Range("I79").Value = 6419.8
$ws.Range("K79").Value = 6419.8
$ws.Range("M79").Value = -5327.8

$ws.Range("H107").Value = 333
$ws.Range("I107").Value = 333
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 333
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1587
$ws.Range("N107").Value = $null

$ws.Range("H132").Value = 73138.516
$ws.Range("I132").Value = 79889.336
$ws.Range("K132").Value = 239668.008
$ws.Range("M132").Value = -237138.008

$ws.Range("H137").Value = 1732985.2
$ws.Range("I137").Value = 1264237.2
$ws.Range("K137").Value = 3792711.6
$ws.Range("M137").Value = -3790161.6

$ws.Range("H138").Value = 2463.2415
$ws.Range("I138").Value = 1929.2273
$ws.Range("K138").Value = 5787.6819
$ws.Range("M138").Value = -647.6818999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1795.2858
$ws.Range("I2").Value = 1613.7778
$ws.Range("J2").Value = 3837.25
$ws.Range("K2").Value = 1613.7778
$ws.Range("L2").Value = 3837.25
$ws.Range("M2").Value = -1500.7778
$ws.Range("N2").Value = -4063.25

$ws.Range("H32").Value = 5498803.5
$ws.Range("I32").Value = 6581699.5
$ws.Range("K32").Value = 6581699.5
$ws.Range("M32").Value = -6581412.5

$ws.Range("H110").Value = 775.4545000000001
$ws.Range("I110").Value = 735.3333
$ws.Range("K110").Value = 735.3333
$ws.Range("M110").Value = 1309.6667

$ws.Range("H116").Value = 1795.2858
$ws.Range("I116").Value = 1613.7778
$ws.Range("J116").Value = 3837.25
$ws.Range("K116").Value = 1613.7778
$ws.Range("L116").Value = 3837.25
$ws.Range("M116").Value = 680.2221999999999
$ws.Range("N116").Value = -8425.25

$ws.Range("H132").Value = 1373475
$ws.Range("J132").Value = 14966.333
$ws.Range("L132").Value = 44898.999
$ws.Range("N132").Value = -49958.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1795.2858
$ws.Range("I3").Value = 1613.7778
$ws.Range("J3").Value = 3837.25
$ws.Range("K3").Value = 1613.7778
$ws.Range("L3").Value = 3837.25
$ws.Range("M3").Value = -1499.7778
$ws.Range("N3").Value = -4065.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6448.2646
$ws.Range("I31").Value = 1773
$ws.Range("K31").Value = 1773
$ws.Range("M31").Value = -1478

$ws.Range("H34").Value = 6448.2646
$ws.Range("I34").Value = 1773
$ws.Range("K34").Value = 1773
$ws.Range("M34").Value = -1571

$ws.Range("H35").Value = 18538.076
$ws.Range("I35").Value = 15999.167
$ws.Range("K35").Value = 15999.167
$ws.Range("M35").Value = -15705.167

$ws.Range("H105").Value = 43751.625
$ws.Range("I105").Value = 49787.57
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 49787.57
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = -48040.57
$ws.Range("N105").Value = -4994

$ws.Range("H107").Value = 727.9375
$ws.Range("I107").Value = 753.6923
$ws.Range("K107").Value = 753.6923
$ws.Range("M107").Value = 1166.3077

$ws.Range("H122").Value = 2971.7273
$ws.Range("I122").Value = 1173.8572
$ws.Range("K122").Value = 3521.5716
$ws.Range("M122").Value = -1071.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 359823.66
$ws.Range("I14").Value = 359823.66
$ws.Range("K14").Value = 1079470.98
$ws.Range("M14").Value = -1079297.98

$ws.Range("H113").Value = 1497.3
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1497.3
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 4491.9
$ws.Range("M113").Value = $null
$ws.Range("N113").Value = -8831.9

$ws.Range("H131").Value = 8406.529
$ws.Range("I131").Value = 1484.1666
$ws.Range("J131").Value = 9889.893
$ws.Range("K131").Value = 4452.4998
$ws.Range("L131").Value = 29669.679
$ws.Range("M131").Value = 587.5002000000004
$ws.Range("N131").Value = -39749.679

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3159.0667
$ws.Range("I113").Value = 1410.5714
$ws.Range("J113").Value = 4689
$ws.Range("K113").Value = 1410.5714
$ws.Range("L113").Value = 4689
$ws.Range("M113").Value = 759.4286
$ws.Range("N113").Value = -9029

$ws.Range("H122").Value = 2875.6667
$ws.Range("I122").Value = 3792.3076
$ws.Range("K122").Value = 11376.9228
$ws.Range("M122").Value = -8926.9228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 36507.83
$ws.Range("J22").Value = 2425.8823
$ws.Range("L22").Value = 2425.8823
$ws.Range("N22").Value = -3015.8823

$ws.Range("H27").Value = 36507.83
$ws.Range("J27").Value = 2425.8823
$ws.Range("L27").Value = 2425.8823
$ws.Range("N27").Value = -2639.8823

$ws.Range("H40").Value = 4850.852
$ws.Range("I40").Value = 4278.96
$ws.Range("K40").Value = 4278.96
$ws.Range("M40").Value = -4142.96

$ws.Range("H46").Value = 2633.2727
$ws.Range("I46").Value = 1357.1428
$ws.Range("J46").Value = 2976.8462
$ws.Range("K46").Value = 1357.1428
$ws.Range("L46").Value = 2976.8462
$ws.Range("M46").Value = -1169.1428
$ws.Range("N46").Value = -3352.8462

$ws.Range("H122").Value = 3110.0635
$ws.Range("I122").Value = 2876.75
$ws.Range("K122").Value = 8630.25
$ws.Range("M122").Value = -6180.25

$ws.Range("H132").Value = 826467.4399999999
$ws.Range("I132").Value = 1050314.6
$ws.Range("J132").Value = 5694.3335
$ws.Range("K132").Value = 3150943.8
$ws.Range("L132").Value = 17083.0005
$ws.Range("M132").Value = -3148413.8
$ws.Range("N132").Value = -22143.0005

$ws.Range("H136").Value = 5558.7144
$ws.Range("I136").Value = 4880.737
$ws.Range("K136").Value = 14642.211
$ws.Range("M136").Value = -12092.211

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2715.3809
$ws.Range("I107").Value = 1040.4615
$ws.Range("J107").Value = 5437.125
$ws.Range("K107").Value = 3121.3845
$ws.Range("L107").Value = 16311.375
$ws.Range("M107").Value = -1201.3845
$ws.Range("N107").Value = -20151.375

$ws.Range("H122").Value = 2774.818
$ws.Range("I122").Value = 2218.7693
$ws.Range("K122").Value = 6656.3079
$ws.Range("M122").Value = -4206.3079
